$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> new value for column C (nombre_aides) and column E (montant_total)
$updates = @{
    2   = @{ C = 100827; E = 327352874 }
    9   = @{ C = 285;    E = 36555237 }
    19  = @{ C = 4367;   E = 66092282 }
    48  = @{ C = 1677;   E = 31738250 }
    56  = @{ C = 11976;  E = 187859672 }
    64  = @{ C = 5212;   E = 20422343 }
    92  = @{ C = 409189; E = 1595779968 }
    93  = @{ C = 209615; E = 1309459328 }
    94  = @{ C = 94218;  E = 918462438 }
    95  = @{ C = 50782;  E = 933328284 }
    97  = @{ C = 2162;   E = 214351518 }
    104 = @{ C = 135253; E = 272253470 }
    119 = @{ C = 356;    E = 10745878 }
    173 = @{ C = 96860;  E = 327935644 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
}

$wb.Save()
